# Auto-generated Excel COM-interop script to apply the Zalera_Profits edit.
# Updates recalculated price/profit figures across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 7115.5386
$ws.Range("I64").Value = 7111.3335
$ws.Range("K64").Value = 7111.3335
$ws.Range("M64").Value = -6863.3335
$ws.Range("H67").Value = 7115.5386
$ws.Range("I67").Value = 7111.3335
$ws.Range("K67").Value = 7111.3335
$ws.Range("M67").Value = -6253.3335
$ws.Range("H69").Value = 2965.6667
$ws.Range("I69").Value = 1885.7142
$ws.Range("J69").Value = 4477.6
$ws.Range("K69").Value = 5657.142599999999
$ws.Range("L69").Value = 13432.8
$ws.Range("M69").Value = -4783.142599999999
$ws.Range("N69").Value = -15180.8
$ws.Range("H72").Value = 2965.6667
$ws.Range("I72").Value = 1885.7142
$ws.Range("J72").Value = 4477.6
$ws.Range("K72").Value = 16971.4278
$ws.Range("L72").Value = 40298.4
$ws.Range("M72").Value = -12603.4278
$ws.Range("N72").Value = -49034.4
$ws.Range("H76").Value = 12504556
$ws.Range("I76").Value = 33336666
$ws.Range("J76").Value = 5289.6
$ws.Range("K76").Value = 33336666
$ws.Range("L76").Value = 5289.6
$ws.Range("M76").Value = -33336351
$ws.Range("N76").Value = -5919.6
$ws.Range("H79").Value = 12504556
$ws.Range("I79").Value = 33336666
$ws.Range("J79").Value = 5289.6
$ws.Range("K79").Value = 33336666
$ws.Range("L79").Value = 5289.6
$ws.Range("M79").Value = -33335574
$ws.Range("N79").Value = -7473.6
$ws.Range("H81").Value = 129989.25
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("H84").Value = 129989.25
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("H86").Value = 4900.6
$ws.Range("I86").Value = 4250.75
$ws.Range("J86").Value = 7500
$ws.Range("K86").Value = 4250.75
$ws.Range("L86").Value = 7500
$ws.Range("M86").Value = -3127.75
$ws.Range("N86").Value = -9746
$ws.Range("H89").Value = 4900.6
$ws.Range("I89").Value = 4250.75
$ws.Range("J89").Value = 7500
$ws.Range("K89").Value = 21253.75
$ws.Range("L89").Value = 37500
$ws.Range("M89").Value = -15637.75
$ws.Range("N89").Value = -48732
$ws.Range("H131").Value = 1843.3572
$ws.Range("I131").Value = 1941.8334
$ws.Range("J131").Value = 1252.5
$ws.Range("K131").Value = 5825.5002
$ws.Range("L131").Value = 3757.5
$ws.Range("M131").Value = -785.5002000000004
$ws.Range("N131").Value = -13837.5
$ws.Range("M81").ClearContents()
$ws.Range("M84").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24889.98
$ws.Range("I32").Value = 24685.6
$ws.Range("K32").Value = 24685.6
$ws.Range("M32").Value = -24398.6
$ws.Range("H61").Value = 6196.6963
$ws.Range("I61").Value = 4866.674
$ws.Range("J61").Value = 12314.8
$ws.Range("K61").Value = 4866.674
$ws.Range("L61").Value = 12314.8
$ws.Range("M61").Value = -4654.674
$ws.Range("N61").Value = -12738.8
$ws.Range("H122").Value = 3323.8125
$ws.Range("I122").Value = 2365
$ws.Range("K122").Value = 7095
$ws.Range("M122").Value = -4645
$ws.Range("H123").Value = 65000
$ws.Range("J123").Value = 65000
$ws.Range("L123").Value = 65000
$ws.Range("N123").Value = -74800
$ws.Range("H132").Value = 4032.4443
$ws.Range("I132").Value = 3532.7908
$ws.Range("J132").Value = 14775
$ws.Range("K132").Value = 10598.3724
$ws.Range("L132").Value = 44325
$ws.Range("M132").Value = -8068.3724
$ws.Range("N132").Value = -49385
$ws.Range("H136").Value = 6196.6963
$ws.Range("I136").Value = 4866.674
$ws.Range("J136").Value = 12314.8
$ws.Range("K136").Value = 14600.022
$ws.Range("L136").Value = 36944.39999999999
$ws.Range("M136").Value = -12050.022
$ws.Range("N136").Value = -42044.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3498.9246
$ws.Range("I134").Value = 2432.4358
$ws.Range("K134").Value = 7297.307400000001
$ws.Range("M134").Value = -4762.307400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 50005184
$ws.Range("J31").Value = 6391.357
$ws.Range("L31").Value = 6391.357
$ws.Range("N31").Value = -6981.357
$ws.Range("H34").Value = 50005184
$ws.Range("J34").Value = 6391.357
$ws.Range("L34").Value = 6391.357
$ws.Range("N34").Value = -6795.357
$ws.Range("H58").Value = 5564.091
$ws.Range("I58").Value = 3751.3333
$ws.Range("J58").Value = 7739.4
$ws.Range("K58").Value = 3751.3333
$ws.Range("L58").Value = 7739.4
$ws.Range("M58").Value = -3548.3333
$ws.Range("N58").Value = -8145.4
$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 498
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -530
$ws.Range("H132").Value = 4559.4546
$ws.Range("I132").Value = 1956.75
$ws.Range("K132").Value = 5870.25
$ws.Range("M132").Value = -3340.25
$ws.Range("H136").Value = 5564.091
$ws.Range("I136").Value = 3751.3333
$ws.Range("J136").Value = 7739.4
$ws.Range("K136").Value = 11253.9999
$ws.Range("L136").Value = 23218.2
$ws.Range("M136").Value = -8703.999899999999
$ws.Range("N136").Value = -28318.2
$ws.Range("N99").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 75508.84
$ws.Range("J37").Value = 75508.84
$ws.Range("L37").Value = 226526.52
$ws.Range("N37").Value = -226750.52
$ws.Range("H93").Value = 3298
$ws.Range("J93").Value = 4247.25
$ws.Range("L93").Value = 12741.75
$ws.Range("N93").Value = -16485.75
$ws.Range("H134").Value = 2087.6667
$ws.Range("I134").Value = 674.7143
$ws.Range("J134").Value = 7033
$ws.Range("K134").Value = 2024.1429
$ws.Range("L134").Value = 21099
$ws.Range("M134").Value = 3045.8571
$ws.Range("N134").Value = -31239
$ws.Range("H139").Value = 55558404
$ws.Range("I139").Value = 71430790
$ws.Range("K139").Value = 214292370
$ws.Range("M139").Value = -214287230
$ws.Range("H140").Value = 1258.0588
$ws.Range("I140").Value = 869
$ws.Range("J140").Value = 1695.75
$ws.Range("K140").Value = 2607
$ws.Range("L140").Value = 5087.25
$ws.Range("M140").Value = 2573
$ws.Range("N140").Value = -15447.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 10718.454
$ws.Range("J19").Value = 11978.111
$ws.Range("L19").Value = 11978.111
$ws.Range("N19").Value = -12554.111
$ws.Range("H102").Value = 2627.2727
$ws.Range("I102").Value = 2490
$ws.Range("K102").Value = 2490
$ws.Range("M102").Value = -868
$ws.Range("H122").Value = 2697
$ws.Range("I122").Value = 2378.8333
$ws.Range("J122").Value = 3333.3333
$ws.Range("K122").Value = 7136.499899999999
$ws.Range("L122").Value = 9999.999899999999
$ws.Range("M122").Value = -4686.499899999999
$ws.Range("N122").Value = -14899.9999
$ws.Range("H132").Value = 5621.7715
$ws.Range("I132").Value = 4331.04
$ws.Range("J132").Value = 8848.6
$ws.Range("K132").Value = 12993.12
$ws.Range("L132").Value = 26545.8
$ws.Range("M132").Value = -10463.12
$ws.Range("N132").Value = -31605.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4915.2856
$ws.Range("I7").Value = 3874.75
$ws.Range("K7").Value = 3874.75
$ws.Range("M7").Value = -3762.75
$ws.Range("H16").Value = 957.5
$ws.Range("I16").Value = 957.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 957.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -787.5
$ws.Range("H126").Value = 4915.2856
$ws.Range("I126").Value = 3874.75
$ws.Range("K126").Value = 11624.25
$ws.Range("M126").Value = -9154.25
$ws.Range("H136").Value = 4954.778
$ws.Range("I136").Value = 3498.923
$ws.Range("K136").Value = 10496.769
$ws.Range("M136").Value = -7946.769
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1375.25
$ws.Range("I62").Value = 1333.6666
$ws.Range("J62").Value = 1500
$ws.Range("K62").Value = 1333.6666
$ws.Range("L62").Value = 1500
$ws.Range("M62").Value = -709.6666
$ws.Range("N62").Value = -2748
$ws.Range("H65").Value = 1375.25
$ws.Range("I65").Value = 1333.6666
$ws.Range("J65").Value = 1500
$ws.Range("K65").Value = 6668.333000000001
$ws.Range("L65").Value = 7500
$ws.Range("M65").Value = -3548.333000000001
$ws.Range("N65").Value = -13740
$ws.Range("H81").Value = 12879.7705
$ws.Range("I81").Value = 2248.4614
$ws.Range("J81").Value = 16828.543
$ws.Range("K81").Value = 4496.9228
$ws.Range("L81").Value = 33657.086
$ws.Range("M81").Value = -3435.9228
$ws.Range("N81").Value = -35779.086
$ws.Range("H84").Value = 12879.7705
$ws.Range("I84").Value = 2248.4614
$ws.Range("J84").Value = 16828.543
$ws.Range("K84").Value = 22484.614
$ws.Range("L84").Value = 168285.43
$ws.Range("M84").Value = -17180.614
$ws.Range("N84").Value = -178893.43
$ws.Range("H126").Value = 3570.647
$ws.Range("I126").Value = 3360.68
$ws.Range("K126").Value = 10082.04
$ws.Range("M126").Value = -7612.039999999999
$ws.Range("H136").Value = 2782.25
$ws.Range("I136").Value = 1242.8182
$ws.Range("K136").Value = 3728.4546
$ws.Range("M136").Value = -1178.4546

